$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.016.63"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "2.505.94"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").Value = "'534.19"
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").Value = "'134.26"
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "'0.571"
$ws.Range("E8").Value = "  +2.68%  "
$ws.Range("D9").Value = "2.511.28"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").Value = "'0.0995"
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("D11").Value = "'0.153"
$ws.Range("E11").Value = "  -2.76%  "
$ws.Range("D12").Value = "'5.19"
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("D13").Value = "'0.331"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "2.951.39"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "58.810.92"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").Value = "'22.38"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "2.507.25"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "'10.63"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("D20").Value = "'4.26"
$ws.Range("E20").Value = "  +1.93%  "
$ws.Range("D21").Value = "'321.60"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("D22").Value = "'6.18"
$ws.Range("E22").Value = "  +1.82%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "'66.05"
$ws.Range("E24").Value = "  +4.24%  "
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D26").Value = "'0.995"
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("D28").Value = "'7.47"
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("D29").Value = "0.0₃0759"
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("D30").Value = "'172.73"
$ws.Range("E30").Value = "  +2.38%  "
$ws.Range("D31").Value = "'1.74"
$ws.Range("E31").Value = "  +1.60%  "
$ws.Range("D32").Value = "'6.29"
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D36").Value = "'18.12"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("E37").Value = "  -3.93%  "
$ws.Range("D38").Value = "'3.97"
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").Value = "'1.52"
$ws.Range("E39").Value = "  +3.32%  "
$ws.Range("D40").Value = "'0.832"
$ws.Range("E40").Value = "  +6.25%  "
$ws.Range("E41").Value = "  -1.42%  "
$ws.Range("D42").Value = "'3.48"
$ws.Range("E42").Value = "  +1.05%  "
$ws.Range("D43").Value = "'275.86"
$ws.Range("E43").Value = "  -2.43%  "
$ws.Range("D44").Value = "'131.92"
$ws.Range("E44").Value = "  +6.30%  "
$ws.Range("D45").Value = "'5.04"
$ws.Range("E45").Value = "  -2.28%  "
$ws.Range("D46").Value = "'0.595"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").Value = "'0.0934"
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("D48").Value = "'0.0511"
$ws.Range("E48").Value = "  +2.52%  "
$ws.Range("D49").Value = "'0.0218"
$ws.Range("E49").Value = "  +2.21%  "
$ws.Range("D50").Value = "'16.80"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("D51").Value = "1.755.97"
$ws.Range("E51").Value = "  +0.64%  "
